# Comparing all four models
# 1) Update existing posterior-summary / model-fit values that changed
# 2) Add four new convergence-diagnostics sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: style a cell like the workbook's existing bold/centered/bordered
# header-and-label style (style index 1 in the original workbook).
# ---------------------------------------------------------------------------
function Set-LabelStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous (renders as "thin")
}

# Helper: set a value into a cell that must remain TEXT even though it looks
# like a number (mirrors the ModelFit_Table sheet, whose numeric-looking
# metrics are stored as inline strings, not numbers).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# ===========================================================================
# PostSummary_TRI_M1 updates
# ===========================================================================
$wsTriM1 = $wb.Worksheets.Item("PostSummary_TRI_M1")
$wsTriM1.Range("B2").Value = -3.75
$wsTriM1.Range("C2").Value = -3.5
$wsTriM1.Range("D2").Value = -3.24

$wsTriM1.Range("B3").Value = -4.18
$wsTriM1.Range("C3").Value = -3.64
$wsTriM1.Range("D3").Value = -3.18

$wsTriM1.Range("B5").Value = 1.01
$wsTriM1.Range("C5").Value = 1.36
$wsTriM1.Range("D5").Value = 1.74

$wsTriM1.Range("B6").Value = -0.2
$wsTriM1.Range("C6").Value = 0.25
$wsTriM1.Range("D6").Value = 1.16

$wsTriM1.Range("B8").Value = 1.02
$wsTriM1.Range("C8").Value = 3.07
$wsTriM1.Range("D8").Value = 6.4

$wsTriM1.Range("B9").Value = -0.04

$wsTriM1.Range("D10").Value = 0.11

# ===========================================================================
# PostSummary_TRI_M2 updates
# ===========================================================================
$wsTriM2 = $wb.Worksheets.Item("PostSummary_TRI_M2")
$wsTriM2.Range("B2").Value = -3.82
$wsTriM2.Range("C2").Value = -3.64
$wsTriM2.Range("D2").Value = -3.46

$wsTriM2.Range("B3").Value = -0.23
$wsTriM2.Range("C3").Value = -0.1
$wsTriM2.Range("D3").Value = 0.04

$wsTriM2.Range("B4").Value = -0.06
$wsTriM2.Range("C4").Value = 0.09
$wsTriM2.Range("D4").Value = 0.23

$wsTriM2.Range("B5").Value = -4.29
$wsTriM2.Range("C5").Value = -3.91
$wsTriM2.Range("D5").Value = -3.58

$wsTriM2.Range("B6").Value = -0.07000000000000001
$wsTriM2.Range("C6").Value = 0.1
$wsTriM2.Range("D6").Value = 0.3

$wsTriM2.Range("B7").Value = 0.03
$wsTriM2.Range("C7").Value = 0.4
$wsTriM2.Range("D7").Value = 0.76

$wsTriM2.Range("D8").Value = 3.39

$wsTriM2.Range("B9").Value = -0.22
$wsTriM2.Range("D9").Value = 0.24

$wsTriM2.Range("B11").Value = 1.17
$wsTriM2.Range("C11").Value = 1.4
$wsTriM2.Range("D11").Value = 1.69

$wsTriM2.Range("B12").Value = -0.13
$wsTriM2.Range("D12").Value = 0.72

$wsTriM2.Range("B14").Value = 1.01
$wsTriM2.Range("C14").Value = 1.94
$wsTriM2.Range("D14").Value = 3.89

$wsTriM2.Range("B15").Value = -0.03

$wsTriM2.Range("D16").Value = 0.11

# ===========================================================================
# ModelFit_Table updates (values stored as text, not numbers)
# ===========================================================================
$wsFit = $wb.Worksheets.Item("ModelFit_Table")
Set-TextValue $wsFit.Range("D4") "0.99"
Set-TextValue $wsFit.Range("E4") "0.97"

Set-TextValue $wsFit.Range("E6") "2.65"

Set-TextValue $wsFit.Range("D7") "0.74"
Set-TextValue $wsFit.Range("E7") "0.80"

Set-TextValue $wsFit.Range("D11") "10.75"
Set-TextValue $wsFit.Range("E11") "11.90"

Set-TextValue $wsFit.Range("D12") "6.96"
Set-TextValue $wsFit.Range("E12") "6.86"

Set-TextValue $wsFit.Range("D13") "6.57"
Set-TextValue $wsFit.Range("E13") "6.77"

# ===========================================================================
# Add four new convergence-diagnostics sheets at the end of the workbook
# ===========================================================================
$convergenceHeaders = @("mean", "sd", "hdi_3%", "hdi_97%", "mcse_mean", "mcse_sd", "ess_bulk", "ess_tail", "r_hat")

function Add-ConvergenceSheet($name, $rows) {
    $wbLocal = $excel.ActiveWorkbook
    $lastSheet = $wbLocal.Worksheets.Item($wbLocal.Worksheets.Count)
    $ws = $wbLocal.Worksheets.Add($null, $lastSheet)
    $ws.Name = $name

    for ($i = 0; $i -lt $convergenceHeaders.Count; $i++) {
        $cell = $ws.Cells.Item(1, $i + 2)
        $cell.Value = $convergenceHeaders[$i]
        Set-LabelStyle $cell
    }

    for ($r = 0; $r -lt $rows.Count; $r++) {
        $rowData = $rows[$r]
        $excelRow = $r + 2

        $labelCell = $ws.Cells.Item($excelRow, 1)
        $labelCell.Value = $rowData[0]
        Set-LabelStyle $labelCell

        for ($c = 1; $c -lt $rowData.Count; $c++) {
            $ws.Cells.Item($excelRow, $c + 1).Value = $rowData[$c]
        }
    }

    $ws.PageSetup.LeftMargin = 54
    $ws.PageSetup.RightMargin = 54
    $ws.PageSetup.TopMargin = 72
    $ws.PageSetup.BottomMargin = 72
    $ws.PageSetup.HeaderMargin = 36
    $ws.PageSetup.FooterMargin = 36
}

# --- Bi_Convergence_M1 -----------------------------------------------------
$biM1Rows = @(
    @("level_2[log_lambda (intercept)]", -3.5444, 0.1033, -3.7341, -3.3423, 0.014, 0.0055, 55.0324, 169.1281, 1.0379),
    @("level_2[log_mu (intercept)]", -3.6375, 0.1906, -4.0056, -3.2987, 0.0209, 0.0124, 82.32940000000001, 92.0359, 1.0553),
    @("level_2[var_log_lambda]", 1.3788, 0.151, 1.1081, 1.6683, 0.0151, 0.007, 101.3883, 266.408, 1.024),
    @("level_2[cov_log_lambda_mu]", 0.217, 0.2602, -0.2543, 0.6837, 0.0399, 0.025, 49.3175, 92.5549, 1.0754),
    @("level_2[var_log_mu]", 2.9116, 1.1705, 0.98, 5.0425, 0.2164, 0.191, 28.2635, 28.4597, 1.1032)
)
Add-ConvergenceSheet "Bi_Convergence_M1" $biM1Rows

# --- Bi_Convergence_M2 -----------------------------------------------------
$biM2Rows = @(
    @("level_2[log_lambda (intercept)]", -3.64, 0.1254, -3.8672, -3.3988, 0.0226, 0.0044, 30.9104, 244.0368, 1.0877),
    @("level_2[log_mu (intercept)]", 0.2056, 0.057, 0.091, 0.3069, 0.0031, 0.0013, 348.7436, 710.9018, 1.0133),
    @("level_2[beta_lambda[0]]", -0.1051, 0.0743, -0.2497, 0.0303, 0.0023, 0.0016, 1063.5685, 1859.6704, 1.0074),
    @("level_2[beta_mu[0]]", 0.08400000000000001, 0.0927, -0.0887, 0.2634, 0.0043, 0.0021, 461.2482, 1018.5832, 1.0052),
    @("level_2[beta_lambda[1]]", -3.9504, 0.2648, -4.4472, -3.4679, 0.0381, 0.0154, 48.8896, 186.3261, 1.0679),
    @("level_2[beta_mu[1]]", 0.0486, 0.1299, -0.2, 0.2886, 0.007900000000000001, 0.0052, 268.9322, 402.8904, 1.0221),
    @("level_2[beta_lambda[2]]", 0.09669999999999999, 0.1068, -0.1011, 0.2961, 0.0059, 0.003, 322.6393, 629.7186, 1.0217),
    @("level_2[beta_mu[2]]", 0.4295, 0.2289, 0.0081, 0.8832, 0.0141, 0.008, 265.7637, 546.0699, 1.0103),
    @("level_2[var_log_lambda]", 1.3803, 0.1672, 1.0688, 1.6917, 0.029, 0.0068, 32.6989, 205.9711, 1.0876),
    @("level_2[cov_log_lambda_mu]", 0.2055, 0.2367, -0.2348, 0.6422, 0.0303, 0.0149, 57.7232, 136.898, 1.0641),
    @("level_2[var_log_mu]", 2.3503, 1.0039, 0.668, 4.3239, 0.1511, 0.0717, 40.4364, 97.0141, 1.0794)
)
Add-ConvergenceSheet "Bi_Convergence_M2" $biM2Rows

# --- Tri_Convergence_M1 -----------------------------------------------------
$triM1Rows = @(
    @("level_2[log_lambda (intercept)]", -3.5331, 0.1106, -3.7339, -3.3262, 0.0209, 0.0067, 28.4573, 115.1459, 1.1),
    @("level_2[log_mu (intercept)]", -3.6239, 0.1963, -3.9877, -3.2549, 0.0219, 0.0143, 81.2829, 140.7322, 1.0492),
    @("level_2[log_eta (intercept)]", 3.2314, 0.0149, 3.2022, 3.2582, 0.0004, 0.0002, 1499.9869, 3149.9149, 1.0011),
    @("level_2[var_log_lambda]", 1.3561, 0.1565, 1.0612, 1.6539, 0.0152, 0.0074, 104.6756, 192.2475, 1.0205),
    @("level_2[cov_log_lambda_mu]", 0.2242, 0.2678, -0.2462, 0.7141, 0.0506, 0.04, 38.5095, 34.83, 1.0935),
    @("level_2[cov_log_lambda_eta]", 0.0131, 0.009900000000000001, -0.0054, 0.0318, 0.0001, 0.0001, 6127.531, 12502.7976, 1.0014),
    @("level_2[var_log_mu]", 2.8646, 1.2543, 1.0011, 5.367, 0.3966, 0.1627, 11.0432, 66.8747, 1.2849),
    @("level_2[cov_log_mu_eta]", -0.008200000000000001, 0.0147, -0.0368, 0.0188, 0.0004, 0.0009, 1379.3299, 636.0331, 1.0141),
    @("level_2[var_log_eta]", 0.08450000000000001, 0.0098, 0.0667, 0.1033, 0.0004, 0.0002, 510.1431, 1001.5409, 1.0083)
)
Add-ConvergenceSheet "Tri_Convergence_M1" $triM1Rows

# --- Tri_Convergence_M2 -----------------------------------------------------
$triM2Rows = @(
    @("level_2[log_lambda (intercept)]", -3.6485, 0.1213, -3.8855, -3.4181, 0.0128, 0.0083, 92.5595, 71.0873, 1.0542),
    @("level_2[log_mu (intercept)]", -0.1, 0.07149999999999999, -0.2403, 0.0297, 0.0028, 0.0022, 685.2626, 855.3287, 1.0067),
    @("level_2[log_eta (intercept)]", 0.0898, 0.08160000000000001, -0.0643, 0.2443, 0.0045, 0.002, 329.423, 611.2028, 1.0045),
    @("level_2[var_log_lambda]", 1.4108, 0.1682, 1.1147, 1.7525, 0.0164, 0.0104, 105.5862, 155.579, 1.0353),
    @("level_2[cov_log_lambda_mu]", 0.1998, 0.2445, -0.2363, 0.6497000000000001, 0.0261, 0.029, 95.3501, 97.31440000000001, 1.0352),
    @("level_2[cov_log_lambda_eta]", 0.0152, 0.0105, -0.0047, 0.0347, 0.0001, 0.0001, 9819.058300000001, 9666.892900000001, 1.001),
    @("level_2[var_log_mu]", 2.0188, 0.8012, 0.8609, 3.5449, 0.172, 0.07389999999999999, 21.7429, 93.524, 1.1361),
    @("level_2[cov_log_mu_eta]", -0.0062, 0.0131, -0.0306, 0.0188, 0.0004, 0.0004, 1141.9562, 1309.9337, 1.0071),
    @("level_2[var_log_eta]", 0.094, 0.01, 0.0752, 0.1124, 0.0005, 0.0002, 478.6025, 1363.5061, 1.0101)
)
Add-ConvergenceSheet "Tri_Convergence_M2" $triM2Rows

# Restore the original active sheet / selection (first sheet, cell A1) so the
# workbook-level view state matches the source edit (bookViews unchanged).
$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("A1").Select() | Out-Null
